$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose target text looks like a plain decimal number
# as Text, so assigning the literal string does not get silently coerced to a
# float (which would also lose significant trailing zeros, e.g. "577.40").
# NumberFormat is set per-cell (a multi-area Range only applies to the first
# area in this host), but every cell reuses the same single new style entry.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '63.569.27'
$ws.Range("E2").Value = '  +5.81%  '

$ws.Range("D3").Value = '3.399.79'
$ws.Range("E3").Value = '  +6.33%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '577.40'
$ws.Range("E5").Value = '  +7.60%  '

$ws.Range("D6").Value = '154.49'
$ws.Range("E6").Value = '  +6.21%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.403.05'
$ws.Range("E8").Value = '  +6.18%  '

$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("D10").Value = '7.48'
$ws.Range("E10").Value = '  +2.04%  '

$ws.Range("E11").Value = '  +7.22%  '

$ws.Range("E12").Value = '  +1.56%  '

$ws.Range("D13").Value = '3.980.28'
$ws.Range("E13").Value = '  +6.20%  '

$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("E15").Value = '  +7.15%  '

$ws.Range("E16").Value = '  +5.00%  '

$ws.Range("D17").Value = '63.629.21'
$ws.Range("E17").Value = '  +5.88%  '

$ws.Range("D18").Value = '3.388.52'
$ws.Range("E18").Value = '  +5.09%  '

$ws.Range("D19").Value = '6.38'
$ws.Range("E19").Value = '  +1.67%  '

$ws.Range("D20").Value = '13.97'
$ws.Range("E20").Value = '  +4.85%  '

$ws.Range("D21").Value = '8.46'
$ws.Range("E21").Value = '  +2.93%  '

$ws.Range("D22").Value = '390.47'
$ws.Range("E22").Value = '  +5.51%  '

$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("D24").Value = '0.538'
$ws.Range("E24").Value = '  +2.83%  '

$ws.Range("D25").Value = '71.01'
$ws.Range("E25").Value = '  +2.19%  '

$ws.Range("D26").Value = '9.71'
$ws.Range("E26").Value = '  +12.53%  '

$ws.Range("E27").Value = '  +18.14%  '

$ws.Range("E28").Value = '  +6.19%  '

$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  +7.63%  '

$ws.Range("D31").Value = '6.44'
$ws.Range("E31").Value = '  +4.98%  '

$ws.Range("D32").Value = '23.20'
$ws.Range("E32").Value = '  +3.16%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '5.59'
$ws.Range("E33").Value = '  +5.90%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.32'
$ws.Range("E34").Value = '  +10.29%  '

$ws.Range("E35").Value = '  +2.91%  '

$ws.Range("E36").Value = '  +8.96%  '

$ws.Range("D37").Value = '158.25'
$ws.Range("E37").Value = '  +0.92%  '

$ws.Range("D38").Value = '27.95'
$ws.Range("E38").Value = '  +5.70%  '

$ws.Range("E39").Value = '  +12.59%  '

$ws.Range("D40").Value = '2.919.24'
$ws.Range("E40").Value = '  +2.81%  '

$ws.Range("E41").Value = '  +5.79%  '

$ws.Range("D42").Value = '0.0327'
$ws.Range("E42").Value = '  +5.74%  '

$ws.Range("D43").Value = '0.763'
$ws.Range("E43").Value = '  +6.05%  '

$ws.Range("D44").Value = '41.14'
$ws.Range("E44").Value = '  +3.00%  '

$ws.Range("E45").Value = '  +1.37%  '

$ws.Range("E46").Value = '  +7.49%  '

$ws.Range("D47").Value = '3.444.73'
$ws.Range("E47").Value = '  +6.36%  '

$ws.Range("E48").Value = '  +6.71%  '

$ws.Range("D49").Value = '301.25'
$ws.Range("E49").Value = '  +13.55%  '

$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("E51").Value = '  +2.57%  '
